# Modifica nombre de "Current Steering" a "W-2W"
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 - Title placeholder: "...Current Steering DAC..." ->
#           "...W" + "-2WCurrent " + "Steering DAC..."
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(1)
$titleTextRange = $titleShape.TextFrame.TextRange

$titleHit = $titleTextRange.Find("Current Steering DAC")
$titleStart = $titleHit.Start
$titleHit.Text = "W-2WCurrent Steering DAC"

# Split the freshly-written text into three runs matching the target
# formatting (the start offset is unchanged by the in-place Text edit).
$runW = $titleTextRange.Characters($titleStart, 1)
$runMid = $titleTextRange.Characters($titleStart + 1, 11)
$runEnd = $titleTextRange.Characters($titleStart + 12, 12)

# Force each chunk to be its own run (re-assigning the same text forces
# PowerPoint to split runs at these boundaries without touching the
# neighbouring runs' formatting).
$runW.Text = "W"
$runMid.Text = "-2WCurrent "
$runEnd.Text = "Steering DAC"

# ---------------------------------------------------------------------
# Slide 2 - "CuadroTexto 1": "Current Steering DAC symbol" ->
#           "W-2W Current " + "Steering DAC symbol"
#           and reposition/resize the caption text box.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$capShape2 = $s2.Shapes.Item(1)
$capRange2 = $capShape2.TextFrame.TextRange

$hit2 = $capRange2.Find("Current Steering DAC symbol")
$start2 = $hit2.Start
$hit2.Text = "W-2W Current Steering DAC symbol"

$run2a = $capRange2.Characters($start2, 13)
$run2b = $capRange2.Characters($start2 + 13, 19)
$run2a.Text = "W-2W Current "
$run2b.Text = "Steering DAC symbol"

$capShape2.Left = 290.3362992125984
$capShape2.Top = 32.934725409448816
$capShape2.Width = 361.90000999999995
$capShape2.Height = 29.081259842519685

# ---------------------------------------------------------------------
# Slide 3 - "CuadroTexto 1": "Current Steering DAC " -> "W-2W DAC "
#           and reposition/resize the caption text box.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$capShape3 = $s3.Shapes.Item(1)
$capRange3 = $capShape3.TextFrame.TextRange

$hit3 = $capRange3.Find("Current Steering DAC ")
$hit3.Text = "W-2W DAC "

$capShape3.Left = 359.1588976377953
$capShape3.Top = 25.591260842519684
$capShape3.Width = 264.2051968503937
$capShape3.Height = 29.081259842519685
